# Sprint backlog updated for day 2
# - User stories 10 ("set a price range for my search") and 11 ("set a
#   distance range for my search") are pushed from Sprint 1 to Sprint 2:
#   their Finish date becomes "PUSHED TO SPRINT 2" and their Status moves
#   from Completed to In-progress in the main table (rows 13 & 14).
# - They are then logged as new rows in the Sprint 2 tracking table,
#   pushing the existing Sprint-2 rows down by two.
# - User story 7 ("best match" sort), already pushed earlier, flips from
#   Completed to In-progress too.
# - A "KEY" legend block is added below the Sprint 2 table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Mark user stories 10 & 11 as pushed to Sprint 2 in the main table ---
$ws.Range("E13").Value = "PUSHED TO SPRINT 2"
$ws.Range("I13").Value = "In-progress"

$ws.Range("E14").Value = "PUSHED TO SPRINT 2"
$ws.Range("I14").Value = "In-progress"

# --- 2. User story 7 ("best match") moves on to In-progress ---
$ws.Range("I16").Value = "In-progress"

# --- 3. Make room in the Sprint 2 tracking table for the two new entries ---
$ws.Rows("22:23").Insert()

# --- 4. Populate the two new Sprint 2 tracking rows ---
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "As a customer, I can set a price range for my search"
$ws.Range("D22").Value = "27/1/2020"
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 1
$ws.Range("H22").Formula = "=(F22/G22)"
$ws.Range("I22").Value = "In-progress"

$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "As a customer, I can set a distance range for my search"
$ws.Range("D23").Value = "28/1/2020"
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 1
$ws.Range("H23").Formula = "=(F23/G23)"
$ws.Range("I23").Value = "In-progress"

# --- 5. Keep the ROI formulas for the whole Sprint 2 table consistent ---
$ws.Range("H20:H30").Formula = "=(F20/G20)"

# --- 6. Add the KEY legend under the tables ---
$ws.Range("B32").Value = "KEY"
$ws.Range("B32").Font.Bold = $true
$ws.Range("B32").Font.Underline = $true

$ws.Range("B33").Value = "code - the medical code for the required procedure"
$ws.Range("B34").Value = "procedure - search using key words to find the procedure required"
$ws.Range("B35").Value = "rating - a customer left review/rating "
$ws.Range("B36").Value = "best match - a self made formula comparing price to distance"

# --- 7. Misc bookkeeping Excel would normally do on its own ---
$ws.Range("K16").Select()
$ws.PageSetup.Orientation = 1
